# Limpieza de documento: eliminar separadores (líneas ────), imágenes
# "inline" decorativas y párrafos vacíos (espaciador antes de tablas).
#
# Estrategia: recorrer la colección $d.Paragraphs, identificar los
# párrafos candidatos a eliminar y después borrar sus Range en orden
# inverso (de atrás hacia adelante) para que los índices de los
# párrafos restantes no se vean afectados por las eliminaciones ya
# realizadas.

$d = $word.ActiveDocument

$wdWithInTable = 12
$dashChar = [char]0x2500   # '─' U+2500 BOX DRAWINGS LIGHT HORIZONTAL

$toDelete = New-Object System.Collections.ArrayList

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1

    # 1) Párrafos que contienen una imagen inline (dibujo) — se eliminan
    #    por completo (el párrafo entero, no solo la imagen).
    if ($p.Range.InlineShapes.Count -gt 0) {
        [void]$toDelete.Add($idx)
        continue
    }

    $txt = $p.Range.Text

    # 2) Párrafos separadores formados por la línea de guiones largos.
    if ($txt -like ("*" + $dashChar + "*")) {
        [void]$toDelete.Add($idx)
        continue
    }

    # 3) Párrafos vacíos (sin texto ni imagen) usados como espaciador
    #    justo a continuación de una tabla (w:spacing w:before="40",
    #    es decir 2pt, y fuera de cualquier celda de tabla).
    if ($txt.Length -le 1 -and $p.Range.InlineShapes.Count -eq 0) {
        $inTable = $p.Range.Information($wdWithInTable)
        if ((-not $inTable) -and ($p.Format.SpaceBefore -eq 2)) {
            [void]$toDelete.Add($idx)
            continue
        }
    }
}

# Eliminar de atrás hacia adelante para no invalidar los índices.
for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $n = $toDelete[$i]
    $d.Paragraphs($n).Range.Delete()
}

Write-Output ("Parrafos eliminados: " + $toDelete.Count)
